$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finish filling in the previously-placeholder match at row 62 (São Paulo x Fortaleza)
$ws.Range("C62").Value = "0 - 0"
$ws.Range("D62").Value = "Empate"

# Append the rest of the matchday 7 fixtures
$data = @(
    @("Fluminense", "Sport", "2 - 1", "Fluminense", 7),
    @("Corinthians", "Internacional", "4 - 2", "Corinthians", 7),
    @("Ceará", "Vitória", "1 - 0", "Ceará", 7),
    @("Bahia", "Botafogo", "1 - 0", "Bahia", 7),
    @("Vasco da Gama", "Palmeiras", "0 - 1", "Palmeiras", 7),
    @("Grêmio", "Santos", "1 - 0", "Grêmio", 7),
    @("Cruzeiro", "Flamengo", "2 - 1", "Cruzeiro", 7),
    @("Bragantino", "Mirassol", " - ", "  -  ", 7)
)

$row = 63
foreach ($match in $data) {
    $ws.Cells.Item($row, 1).Value = $match[0]
    $ws.Cells.Item($row, 2).Value = $match[1]
    $ws.Cells.Item($row, 3).Value = $match[2]
    $ws.Cells.Item($row, 4).Value = $match[3]
    $ws.Cells.Item($row, 5).Value = $match[4]
    $row = $row + 1
}
